# Commit via gitrun.py em 2024-09-21 18:00:49
# Update "Inscritos" (E), "Pagos" (F) and "Inscricoes homologadas" (H)
# figures for a handful of rows in the Table1 data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 426

$ws.Range("F11").Value = 162
$ws.Range("H11").Value = 162

$ws.Range("E12").Value = 424
$ws.Range("F12").Value = 227
$ws.Range("H12").Value = 227

$ws.Range("E13").Value = 113

$ws.Range("E14").Value = 105
$ws.Range("F14").Value = 53
$ws.Range("H14").Value = 53

$ws.Range("E16").Value = 176

$ws.Range("E17").Value = 82

$ws.Range("E22").Value = 149

$ws.Range("E24").Value = 183
$ws.Range("F24").Value = 96
$ws.Range("H24").Value = 96

$ws.Range("E25").Value = 228

$ws.Range("E26").Value = 127

$ws.Range("E27").Value = 290
$ws.Range("F27").Value = 138
$ws.Range("H27").Value = 138

$ws.Range("E29").Value = 150

$ws.Range("E30").Value = 188
$ws.Range("F30").Value = 108
$ws.Range("H30").Value = 108

$ws.Range("E33").Value = 254

$ws.Range("E34").Value = 189
$ws.Range("F34").Value = 115
$ws.Range("H34").Value = 115

$ws.Range("E35").Value = 124

$ws.Range("E37").Value = 135

$ws.Range("E41").Value = 346

$ws.Range("E42").Value = 318
$ws.Range("F42").Value = 167
$ws.Range("H42").Value = 167

$ws.Range("E44").Value = 272
$ws.Range("F44").Value = 131
$ws.Range("H44").Value = 131

$ws.Range("E45").Value = 125

$ws.Range("E46").Value = 274
$ws.Range("F46").Value = 149
$ws.Range("H46").Value = 149

$ws.Range("E47").Value = 387

$ws.Range("E49").Value = 258
$ws.Range("F49").Value = 108
$ws.Range("H49").Value = 108

$ws.Range("E50").Value = 223

$ws.Range("E51").Value = 211
$ws.Range("F51").Value = 86
$ws.Range("H51").Value = 86
